{"js": "// The two placeholder pictures in the \"Size\" section are replaced with plain\n// hyperlink runs whose visible text is the image's real URL (the images\n// themselves were only ever tiny 1x1 placeholder pixels, keyed by their\n// alt-text description).\nconst urlByDescription = {\n  \"Width of underground pedestrian walkway with single-loaded uses\":\n    \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C18_Underground_Link_A.jpg?h=100%25&w=100%25\",\n  \"Width of underground pedestrian walkway with double-loaded uses\":\n    \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C19_Underground_Link_B.jpg?h=1896&w=3022\",\n};\n\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\npictures.items.forEach((picture) => {\n  picture.load(\"altTextDescription\");\n});\nawait context.sync();\n\nfor (const picture of pictures.items) {\n  const url = urlByDescription[picture.altTextDescription];\n  if (!url) {\n    continue;\n  }\n  // Replace the picture with a text run containing the URL, then turn that\n  // same range into a hyperlink pointing at the URL (Word applies the\n  // built-in \"Hyperlink\" character style automatically).\n  const range = picture.getRange().insertText(url, \"Replace\");\n  range.hyperlink = url;\n}\n\nawait context.sync();\n", "ps1": "# The two placeholder pictures in the \"Size\" section are replaced with plain\n# hyperlink runs whose visible text is the image's real URL (the images\n# themselves were only ever tiny 1x1 placeholder pixels, keyed by their\n# alt-text description).\n$d = $word.ActiveDocument\n\n$urlByDescription = @{\n    \"Width of underground pedestrian walkway with single-loaded uses\" = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C18_Underground_Link_A.jpg?h=100%25&w=100%25\"\n    \"Width of underground pedestrian walkway with double-loaded uses\" = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C19_Underground_Link_B.jpg?h=1896&w=3022\"\n}\n\n# Walk the live InlineShapes collection. Converting a matched picture into a\n# hyperlink run removes it from InlineShapes, so the collection shrinks as we\n# go; only advance the index when the current shape is left untouched.\n$i = 1\nwhile ($i -le $d.InlineShapes.Count) {\n    $shape = $d.InlineShapes.Item($i)\n    $description = $shape.AlternativeText\n    if ($urlByDescription.ContainsKey($description)) {\n        $url = $urlByDescription[$description]\n        $range = $shape.Range\n        # Replace the picture with a text run containing the URL, then turn\n        # that same range into a hyperlink pointing at the URL (Word applies\n        # the built-in \"Hyperlink\" character style automatically).\n        $range.Text = $url\n        $d.Hyperlinks.Add($range, $url) | Out-Null\n    } else {\n        $i = $i + 1\n    }\n}\n"}
